$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J
$ws.Range("J1").Value = "v_1"

# Flip sign of column I values (rows 2-6)
$ws.Range("I2").Value = -23
$ws.Range("I3").Value = -19
$ws.Range("I4").Value = -11
$ws.Range("I5").Value = -7
$ws.Range("I6").Value = -4

# New column J values (rows 2-6), stored as text
$ws.Range("J2").Value = "909.0929545431534∠0.0002494634968685872"
$ws.Range("J3").Value = "861.0836196328438∠0.00024311240502962285"
$ws.Range("J4").Value = "802.3895562630411∠0.0005436034131171138"
$ws.Range("J5").Value = "797.1030359997383∠0.0006786000893114781"
$ws.Range("J6").Value = "787.0108321999132∠0.0009269517133794074"
